# "customer registration test committed"
#
# The Login sheet's sample credentials are replaced with a single
# customer-registration style login: a real-looking admin e-mail in A1
# (still a mailto hyperlink) and a plain numeric value in B1 (its former
# "password" hyperlink is removed, though the cell keeps hyperlink-like
# colouring without the underline).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop every hyperlink on the sheet (A1's old address + B1's password
# link) and re-create only the one that should survive, now pointing at
# the new address.
$ws.Hyperlinks.Delete()
$ws.Range("A1").Value = "admin@itwinetech.com"
$ws.Hyperlinks.Add($ws.Range("A1"), "mailto:admin@itwinetech.com")

# B1 stops being a hyperlink/string and becomes a plain number; keep the
# hyperlink-colour font but drop the underline now that it's not a link.
$ws.Range("B1").Font.Underline = $false
$ws.Range("B1").Value = 1234

# New print setup for the sheet (A4, portrait).
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
